$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from AC1 (an existing header cell) onto the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# Add new header cells for Wins, Losses, Ties in row 1 (columns AD, AE, AF)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in season record values for all data rows (2 through 44):
# Wins = 71, Losses = 91, Ties = 0 for every team/player row
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value = 71
    $ws.Cells.Item($r, 31).Value = 91
    $ws.Cells.Item($r, 32).Value = 0
}
